# Apply the edits described by the commit diff:
#  - Bus connections!E2 and E3: 7500 -> 10000 (dependent formulas in C2/C3 recalc automatically)
#  - Generator data!E4: 15000 -> 20000 (dependent formulas in F4/G4/H4 recalc automatically)
#  - Active sheet changes from "Bus index" to "Bus connections"
#  - Selection on "Bus connections" becomes E4 (single cell)
#  - Selection on "Generator data" becomes E5 (single cell)

$wb = $excel.ActiveWorkbook

$wsBusConnections = $wb.Worksheets.Item("Bus connections")
$wsGeneratorData   = $wb.Worksheets.Item("Generator data")

# Update the underlying data values.
$wsBusConnections.Range("E2").Value = 10000
$wsBusConnections.Range("E3").Value = 10000
$wsGeneratorData.Range("E4").Value = 20000

# Update selections on sheets that are not becoming the active tab first,
# so the final activation/selection below ends up as the active sheet.
$wsGeneratorData.Range("E5").Select()

# Make "Bus connections" the active sheet and set its selection, matching
# the new tabSelected/activeTab state in the workbook.
$wsBusConnections.Activate()
$wsBusConnections.Range("E4").Select()
